$d = $word.ActiveDocument
$d.Content.Find.Execute("61×33=2013", $true, $false, $false, $false, $false, $true, 1, $false, "32×14=448", 2) | Out-Null
$d.Content.Find.Execute("45×59=2655", $true, $false, $false, $false, $false, $true, 1, $false, "24×90=2160", 2) | Out-Null
$d.Content.Find.Execute("66×61=4026", $true, $false, $false, $false, $false, $true, 1, $false, "19×41=779", 2) | Out-Null
$d.Content.Find.Execute("83×23=1909", $true, $false, $false, $false, $false, $true, 1, $false, "42×47=1974", 2) | Out-Null
$d.Content.Find.Execute("75×32=2400", $true, $false, $false, $false, $false, $true, 1, $false, "80×38=3040", 2) | Out-Null
$d.Content.Find.Execute("22×84=1848", $true, $false, $false, $false, $false, $true, 1, $false, "66×70=4620", 2) | Out-Null
$d.Content.Find.Execute("51×13=663", $true, $false, $false, $false, $false, $true, 1, $false, "76×18=1368", 2) | Out-Null
$d.Content.Find.Execute("23×90=2070", $true, $false, $false, $false, $false, $true, 1, $false, "25×75=1875", 2) | Out-Null
$d.Content.Find.Execute("87×21=1827", $true, $false, $false, $false, $false, $true, 1, $false, "66×95=6270", 2) | Out-Null
$d.Content.Find.Execute("52×57=2964", $true, $false, $false, $false, $false, $true, 1, $false, "82×61=5002", 2) | Out-Null
$d.Content.Find.Execute("82×11=902", $true, $false, $false, $false, $false, $true, 1, $false, "29×87=2523", 2) | Out-Null
$d.Content.Find.Execute("25×59=1475", $true, $false, $false, $false, $false, $true, 1, $false, "57×94=5358", 2) | Out-Null
$d.Content.Find.Execute("47×99=4653", $true, $false, $false, $false, $false, $true, 1, $false, "75×75=5625", 2) | Out-Null
$d.Content.Find.Execute("48×16=768", $true, $false, $false, $false, $false, $true, 1, $false, "20×28=560", 2) | Out-Null
$d.Content.Find.Execute("20×59=1180", $true, $false, $false, $false, $false, $true, 1, $false, "65×81=5265", 2) | Out-Null
$d.Content.Find.Execute("75×87=6525", $true, $false, $false, $false, $false, $true, 1, $false, "18×13=234", 2) | Out-Null
$d.Content.Find.Execute("67×22=1474", $true, $false, $false, $false, $false, $true, 1, $false, "27×62=1674", 2) | Out-Null
$d.Content.Find.Execute("14×70=980", $true, $false, $false, $false, $false, $true, 1, $false, "97×49=4753", 2) | Out-Null
$d.Content.Find.Execute("74×57=4218", $true, $false, $false, $false, $false, $true, 1, $false, "86×67=5762", 2) | Out-Null
$d.Content.Find.Execute("27×23=621", $true, $false, $false, $false, $false, $true, 1, $false, "96×29=2784", 2) | Out-Null
$d.Content.Find.Execute("41×93=3813", $true, $false, $false, $false, $false, $true, 1, $false, "33×11=363", 2) | Out-Null
$d.Content.Find.Execute("96×41=3936", $true, $false, $false, $false, $false, $true, 1, $false, "55×66=3630", 2) | Out-Null
$d.Content.Find.Execute("26×89=2314", $true, $false, $false, $false, $false, $true, 1, $false, "87×28=2436", 2) | Out-Null
$d.Content.Find.Execute("50×37=1850", $true, $false, $false, $false, $false, $true, 1, $false, "97×19=1843", 2) | Out-Null
$d.Content.Find.Execute("85×54=4590", $true, $false, $false, $false, $false, $true, 1, $false, "49×85=4165", 2) | Out-Null
